# forests-scraped.xlsx update - 2025-10-20 12:18
# 1) The 4 rows currently on "New" (rows 2-5) are the *previously* scraped
#    listings - they move down to the bottom of "Previously added" (becoming
#    rows 191-194).
# 2) "New" is repopulated with 5 freshly scraped listings.

$wb = $excel.ActiveWorkbook
$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew  = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------
# Part 1: move the current "New" rows (2-5) to the end of "Previously added"
# ---------------------------------------------------------------------

$lastPrevRow = $wsPrev.Cells.Item(1, 1).Worksheet.UsedRange.Rows.Count
# UsedRange always starts at row 1 here, so the used row count is also the
# last used row number.
$destFirst = $lastPrevRow + 1
$destLast = $destFirst + 3

# Capture the hyperlink addresses belonging to the rows we are about to move
# before we touch anything. (Range-scoped .Hyperlinks enumerations return
# unusable proxies in this host, so walk the sheet-level collection and
# match by the hyperlink's own Range address instead.)
$urlByAddr = @{}
foreach ($h in $wsNew.Hyperlinks) {
    $urlByAddr[$h.Range().Address()] = $h.Address()
}
$urls = @($urlByAddr["`$A`$2"], $urlByAddr["`$A`$3"], $urlByAddr["`$A`$4"], $urlByAddr["`$A`$5"])

# Copy the values (this preserves the shared-string text, including the
# numeric-looking cadastre numbers as text) into the new rows.
$wsNew.Range("A2:F5").Copy()
$wsPrev.Range("A" + $destFirst + ":F" + $destLast).PasteSpecial(-4163)

# Copy the formatting from the last pre-existing row so the new rows look
# identical to the rest of the sheet.
$wsPrev.Range("A" + $lastPrevRow + ":F" + $lastPrevRow).Copy()
$wsPrev.Range("A" + $destFirst + ":F" + $destLast).PasteSpecial(-4122)

# Re-create the hyperlinks on column A of the moved rows.
for ($i = 0; $i -lt 4; $i++) {
    $r = $destFirst + $i
    $cell = $wsPrev.Range("A" + $r)
    $txt = $cell.Value()
    $wsPrev.Hyperlinks.Add($cell, $urls[$i], [Type]::Missing, [Type]::Missing, $txt)
}

# Adding hyperlinks resets the cell style to the built-in "Hyperlink" style,
# so re-apply the original formatting once more.
$wsPrev.Range("A" + $lastPrevRow + ":F" + $lastPrevRow).Copy()
$wsPrev.Range("A" + $destFirst + ":F" + $destLast).PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Part 2: replace the listings on "New" with the freshly scraped ones
# ---------------------------------------------------------------------

# Drop the old hyperlinks + values from rows 2-5 first.
$wsNew.Range("A2:F5").Hyperlinks.Delete()
$wsNew.Range("A2:F5").ClearContents()

$newRows = @(
    @{ link="https://www.ss.com/msg/lv/real-estate/wood/aizkraukle-and-reg/aizkraukles-pag/fppdj.html"; price="20 000 €"; district="Aizkraukle un raj."; area="2 ha.";  cadastre="32780050049"; date=45948.00208333333 },
    @{ link="https://www.ss.com/msg/lv/real-estate/wood/aizkraukle-and-reg/aizkraukles-pag/ljdbo.html"; price="30 000 €"; district="Aizkraukle un raj."; area="2 ha.";  cadastre="32780050049"; date=45948.00208333333 },
    @{ link="https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/mezvidu-pag/exmxd.html";         price="6 000 €";  district="Ludza un raj.";     area="1 ha.";  cadastre="68700050016"; date=45949.69305555556 },
    @{ link="https://www.ss.com/msg/lv/real-estate/wood/ogre-and-reg/ledmanes-pag/jexbj.html";         price="31 000 €"; district="Ogre un raj.";      area="1 ha.";  cadastre="74640020009"; date=45949.60902777778 },
    @{ link="https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/cornajas-pag/afdcl.html";      price="15 000 €"; district="Rēzekne un raj.";    area="2 ha.";  cadastre="78460100069"; date=45947.663888888885 }
)

$row = 2
foreach ($item in $newRows) {
    $rowRef = [string]$row

    # Column order matches how the sheet was actually authored (left to
    # right), so new shared-string entries land in the same order as the
    # reference workbook.
    $aCell = $wsNew.Range("A" + $rowRef)
    $aCell.Value = $item.link
    $wsNew.Hyperlinks.Add($aCell, $item.link)

    $wsNew.Range("B" + $rowRef).Value = $item.price
    $wsNew.Range("C" + $rowRef).Value = $item.district
    $wsNew.Range("D" + $rowRef).Value = $item.area

    # Force the cadastre number to be stored as text, not a number.
    $eCell = $wsNew.Range("E" + $rowRef)
    $eCell.NumberFormat = "@"
    $eCell.Value = $item.cadastre

    $wsNew.Range("F" + $rowRef).Value = $item.date

    # Re-apply the standard row formatting (hyperlink-add + forced text
    # format both clobber styles).
    $wsPrev.Range("A" + $lastPrevRow + ":F" + $lastPrevRow).Copy()
    $wsNew.Range("A" + $rowRef + ":F" + $rowRef).PasteSpecial(-4122)

    $row++
}

# Clean up the unused built-in "Hyperlink" cell style that Hyperlinks.Add
# creates (it is not referenced by any cell once formats are reapplied).
foreach ($s in $wb.Styles) {
    if ($s.Name() -eq "Hyperlink") {
        $s.Delete()
    }
}

$excel.CutCopyMode = 0
